$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DEC-2020")

# --- Row 15: fill in the previously-empty row with data (same look as rows 12/13) ---
# Copy formatting from row 13 (same style pattern: s="2" customFormat="1", no explicit height)
$ws.Range("A13:G13").Copy()
$ws.Range("A15:G15").PasteSpecial(-4122)

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 44152
$ws.Range("C15").Value = " Selenium log files"
$ws.Range("D15").Value = "Selenium log file Testing (QMVAR TO GSPN)"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = "Completed"

# --- Row 16: brand-new row appended after row 15 (same look as row 14, taller row) ---
# Copy formatting from row 14 (style pattern without row-level s="2", taller row height)
$ws.Range("A14:G14").Copy()
$ws.Range("A16:G16").PasteSpecial(-4122)

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 44153
$ws.Range("C16").Value = " Selenium log files, Soniya"
# Set G16 before D16 so new shared-string entries are appended in the same
# order as the target workbook (142 = G16 text, 143 = D16 text).
$ws.Range("G16").Value = "Bic_Report_Soukastu Setup Create"
$ws.Range("D16").Value = "Selenium log file Testing (QMVAR TO GSPN), Bic_Report_Soukastu"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = "Completed"

$ws.Rows.Item(16).RowHeight = 28.8

# --- Update the active selection shown when the sheet is opened ---
$ws.Activate()
$ws.Range("E10").Select() | Out-Null
